$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5:E6").NumberFormat = "@"

$ws.Range("A5").Value = "10034"
$ws.Range("B5").Value = "Equal Exchange - Black Silk Espresso"
$ws.Range("C5").Value = "10"
$ws.Range("D5").Value = "71.50"
$ws.Range("E5").Value = "715.00"

$ws.Range("A6").Value = "10400"
$ws.Range("B6").Value = "Equal Exchange - One World"
$ws.Range("C6").Value = "5"
$ws.Range("D6").Value = "71.50"
$ws.Range("E6").Value = "357.50"

$ws.Range("A5:E6").Style = "Normal"
